# Fix the main bug, the project can reluctantly run through the process.
#
# This script applies the bug-tracking sheet update:
#  - removes the obsolete "为什么client的delayTime..." row (old row 17)
#    which shifts every following row up by one
#  - adds two freshly-triaged bug rows at the end of the table
#    (5.1 surround channel crash + duplicate client/server delay timestamp bug)
#  - refreshes the row heights / selection to match the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# 1. Drop the resolved bug row; Excel shifts rows 18..35 up to 17..34
#    automatically, carrying their values/styles/heights with them.
$ws.Rows.Item(17).Delete()

# 2. New bug #18 - duplicate delay time on both devices.
$ws.Rows.Item(22).RowHeight = 59
$ws.Range("A22").Value = 18
$ws.Range("B22").Value = "两台设备发送延时是相同的。。这个不对"
$ws.Range("C22").Value = "指令时间戳与延时的生成应该放在发送指令的时候，不然就是一样的时间"
$ws.Range("C22").WrapText = $true
$ws.Range("D22").Value = "√"

# 3. New bug #19 - 5.1 surround right channel crash (array out of bounds).
$ws.Rows.Item(23).RowHeight = 43
$ws.Range("B23").Value = "分配5.1立体声的右方声道会崩溃"
$ws.Range("C23").Value = "数组越界了，m_lastId只开到了6也就是目前设置的最大id数，但声道是有8个的，所以这里应与声道数保持一致"
$ws.Range("C23").WrapText = $true
$ws.Range("D23").Value = "√"

# 4. Update the sheet view's current selection to match where the author
#    ended up after triaging the new rows.
$ws.Range("C27").Select()

# Best-effort: keep the stored window size in sync with the authored
# workbook (no-op on engines that don't expose this knob).
$win = $wb.Windows.Item(1)
$win.Width = 16608
$win.Height = 27060
